# Expand the EDA notes:
#  - strip the stray _GoBack bookmark off the "It help understand..." paragraph
#  - turn "Contain some steps like:" + the single bulleted list item into a
#    plain paragraph "Contain some steps like Proportion of missing value"
#  - append new notes: a bold "Understanding the data" heading, a
#    "Summary of the features -> df.info" line, a "Show non null entries and
#    feature type -> df.info" line (carrying the relocated _GoBack bookmark),
#    and a trailing empty paragraph.

$d = $word.ActiveDocument

# Locate the two anchor paragraphs by their current text.
$pStart = $null
$pEnd = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "It help understand the dataset*") { $pStart = $p }
    if ($t -like "Proportion of missing value*") { $pEnd = $p }
}

$range = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml =  "<w:p $wns>"
$xml += "<w:r><w:t>It help understand the dataset and identify the issue that could affect model performance downstream.</w:t></w:r>"
$xml += "<w:r><w:t xml:space=""preserve""> </w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $wns>"
$xml += "<w:r><w:t>Contain some steps like</w:t></w:r>"
$xml += "<w:r><w:t xml:space=""preserve""> </w:t></w:r>"
$xml += "<w:r><w:t>Proportion of missing value</w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $wns>"
$xml += "<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>"
$xml += "<w:r><w:rPr><w:b/></w:rPr><w:t>Understanding the data</w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $wns>"
$xml += "<w:r><w:t>Summary of the features -&gt; df.info</w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $wns>"
$xml += "<w:r><w:t xml:space=""preserve"">Show non null entries and feature type </w:t></w:r>"
$xml += "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/>"
$xml += "<w:r><w:t>-&gt; df.info</w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $wns/>"

$range.InsertXML($xml)
